# "fix: various task transformer improvements" - xlsx sample update.
#
# The sample sheet gains three new header columns ("goal_version",
# "rule_name_id", "rule_version") right after the existing "goal_name_id"
# column. Concretely this means:
#   - a single new (blank) column is inserted just before the existing
#     "Parameter / [optional parameter]" column, which pushes that column
#     and the "Values / default, [alternatives]" column one slot to the
#     right
#   - the three newly available columns get the new header text, using
#     the same look & feel as the neighbouring "goal_name_id" header
#   - the selection/active cell is updated to match the authored file

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single new column right before the "Parameter" header column
# (column AP). Everything at/after that column - including the
# "Parameter" and "Values" header columns - shifts one column to the
# right (AR->AS, AS->AT).
$ws.Columns("AP").Insert()

# Give the three newly available header cells (AO1:AQ1) the same
# formatting as the neighbouring "goal_name_id" header cell (AN1), then
# fill in the new header text.
$ws.Range("AN1").Copy()
$ws.Range("AO1:AQ1").PasteSpecial(-4122)

$ws.Range("AO1").Value2 = "goal_version"
$ws.Range("AP1").Value2 = "rule_name_id"
$ws.Range("AQ1").Value2 = "rule_version"

# The data row (row 2) picks up matching (empty) formatting under the new
# "goal_version" column, just like its neighbours.
$ws.Range("AN2").Copy()
$ws.Range("AO2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Keep the new column's width consistent with its neighbour.
try {
    $ws.Range("AP1").EntireColumn.ColumnWidth = 14.15
} catch {
}

# Match the saved selection/active cell from the authored workbook.
$ws.Range("AQ1").Select() | Out-Null

# Best-effort cosmetic view-state tweaks (not all view properties are
# guaranteed to round-trip through this runtime).
try {
    $excel.ActiveWindow.ScrollColumn = 35
} catch {
}
try {
    $ws.StandardWidth = 11.625
} catch {
}
